$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.135057687759399
$ws.Range("B1").Value = 2.677855491638184
$ws.Range("C1").Value = 3.810995101928711
$ws.Range("D1").Value = 5.543853759765625
$ws.Range("E1").Value = 1.788248538970947
